$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet "Sheet1" -> "arima_graph"
# ---------------------------------------------------------------------------
$wsArima = $wb.Worksheets.Item(1)
$wsArima.Name = "arima_graph"

# Update the embedded chart's series formulas so they reference the renamed
# sheet instead of the old "Sheet1" name.
$chartObj = $wsArima.ChartObjects(1)
$chart = $chartObj.Chart

$ser1 = $chart.SeriesCollection(1)
$ser1.Formula = "=SERIES(arima_graph!`$B`$1,,arima_graph!`$B`$2:`$B`$13,1)"

$ser2 = $chart.SeriesCollection(2)
$ser2.Formula = "=SERIES(arima_graph!`$C`$1,,arima_graph!`$C`$2:`$C`$13,2)"

# ---------------------------------------------------------------------------
# 2. Add the new "pm10_limits" worksheet after "arima_graph"
# ---------------------------------------------------------------------------
$wsLimits = $wb.Worksheets.Add($null, $wsArima)
$wsLimits.Name = "pm10_limits"

# Header row (row 4)
$wsLimits.Range("F4").Value = "ue"
$wsLimits.Range("G4").Value = "oms"

# Row 5
$wsLimits.Range("E5").Value = "Promedio 24 horas"
$wsLimits.Range("F5").Value = "<= 50 ug/m3; <= 35 días"
$wsLimits.Range("G5").Value = "<= 50 ug/m3; <= 3 días"

# Row 6
$wsLimits.Range("E6").Value = "Promedio anual"
$wsLimits.Range("F6").Value = "<= 40 ug/m3"
$wsLimits.Range("G6").Value = "<= 20 ug/m3"

# Column widths
$wsLimits.Range("E5").ColumnWidth = 19.1796875
$wsLimits.Range("F5").ColumnWidth = 20.86328125
$wsLimits.Range("G5").ColumnWidth = 22.1796875

# Center-align the "ue"/"oms" value cells (F4:G6)
$wsLimits.Range("F4:G6").HorizontalAlignment = -4108

# Selection / active cell on the new sheet
$wsLimits.Range("G9").Select()

$wsLimits.Activate()
